{"js": "// Replace the lyrics of the single-paragraph song with the new\n// \"Sajeev / Nasi Goreng\" lyrics, line-by-line, while preserving the\n// original <w:br/> line-break structure inside the run.\n//\n// The whole song lives in ONE paragraph, as a single run containing many\n// <w:t> segments separated by <w:br/>. Because several lines repeat\n// verbatim (e.g. \"Yaseen, Yaseen\" / \"Chorus:\") a plain text search/replace\n// would be ambiguous, so we walk the paragraph's text ranges in document\n// order (splitting on the vertical-tab \\u000b that Office.js uses to\n// represent a <w:br/>) and rewrite each non-blank line by position.\n\nconst oldLines = [\n  \"Verse 1:\",\n  \"Yaseen, my Pakistani friend\",\n  \"With his warm smile that never ends\",\n  \"From Lahore to Karachi, he's a true gem\",\n  \"Always there for me till the very end\",\n  \"Chorus:\",\n  \"Yaseen, Yaseen\",\n  \"A friend so true and keen\",\n  \"In his heart, love is seen\",\n  \"Yaseen, Yaseen\",\n  \"Verse 2:\",\n  \"I remember the days we spent\",\n  \"Laughing and talking 'til the night was spent\",\n  \"His kindness and wisdom, a guiding light\",\n  \"In his presence, everything feels right\",\n  \"Chorus:\",\n  \"Yaseen, Yaseen\",\n  \"A friend so true and keen\",\n  \"In his heart, love is seen\",\n  \"Yaseen, Yaseen\",\n  \"Bridge:\",\n  \"No matter the distance or time apart\",\n  \"Yaseen will always have a place in my heart\",\n  \"With his loyalty and friendship, I am blessed\",\n  \"Forever grateful for his love and zest\",\n  \"Chorus:\",\n  \"Yaseen, Yaseen\",\n  \"A friend so true and keen\",\n  \"In his heart, love is seen\",\n  \"Yaseen, Yaseen\",\n  \"Outro:\",\n  \"To my dear Pakistani friend\",\n  \"Yaseen, may our bond never end\",\n  \"In this world or the next, I'll always be\",\n  \"Grateful for your friendship eternally.\",\n];\n\nconst newLines = [\n  \"(Verse 1)\",\n  \"Sajeev, my Malaysian friend,\",\n  \"With a love for Nasi Goreng that never ends,\",\n  \"He cooks it up with skill and flair,\",\n  \"A taste of home that he loves to share.\",\n  \"(Chorus)\",\n  \"Nasi Goreng, oh Nasi Goreng,\",\n  \"It's Sajeev's favorite thing,\",\n  \"Spicy, savory, and oh so good,\",\n  \"In every bite, you can taste the love.\",\n  \"(Verse 2)\",\n  \"From Kuala Lumpur to Penang street,\",\n  \"Sajeev knows where to find the best to eat,\",\n  \"In hawker stalls and busy markets,\",\n  \"He'll hunt down the dish that steals his heart.\",\n  \"(Chorus)\",\n  \"Nasi Goreng, oh Nasi Goreng,\",\n  \"It's Sajeev's favorite thing,\",\n  \"Spicy, savory, and oh so good,\",\n  \"In every bite, you can taste the love.\",\n  \"(Bridge)\",\n  \"With a side of sambal and a fried egg on top,\",\n  \"Sajeev's Nasi Goreng will make your taste buds pop,\",\n  \"A symphony of flavors in every spoonful,\",\n  \"It's a dish that's truly beautiful.\",\n  \"(Chorus)\",\n  \"Nasi Goreng, oh Nasi Goreng,\",\n  \"It's Sajeev's favorite thing,\",\n  \"Spicy, savory, and oh so good,\",\n  \"In every bite, you can taste the love.\",\n  \"(Outro)\",\n  \"So here's to Sajeev and his Nasi Goreng delight,\",\n  \"A dish that brings joy and a smile so bright,\",\n  \"May his love for it never wane,\",\n  \"And may he forever enjoy the taste of home again and again.\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\n\n// Split the paragraph into text ranges on the vertical-tab (\\u000b), which\n// is how Office.js represents a <w:br/> line break. Each returned chunk\n// keeps its trailing delimiter (the \\u000b, or \\r for the very last one);\n// we only ever rewrite the textual part of each chunk, see below.\nconst lineRanges = paragraph.getTextRanges([\"\\u000b\"], false);\nlineRanges.load(\"text\");\nawait context.sync();\n\nlet lineIndex = 0;\nfor (const range of lineRanges.items) {\n  const raw = range.text;\n  // A chunk either ends with \\u000b (a <w:br/> line break, which we must\n  // preserve) or, for the very last chunk, with \\r (the paragraph mark,\n  // which must NOT be re-inserted or Office.js will split the paragraph).\n  const endsWithBreak = raw.length > 0 && raw[raw.length - 1] === \"\\u000b\";\n  const endsWithParaMark = raw.length > 0 && raw[raw.length - 1] === \"\\r\";\n  const trailer = endsWithBreak ? \"\\u000b\" : \"\";\n  const text = (endsWithBreak || endsWithParaMark) ? raw.slice(0, -1) : raw;\n\n  // Blank lines (the empty line between stanzas) are left untouched.\n  if (text.length === 0) {\n    continue;\n  }\n\n  if (lineIndex >= oldLines.length) {\n    throw new Error(\"More lyric lines found than expected: \" + JSON.stringify(raw));\n  }\n  if (text !== oldLines[lineIndex]) {\n    throw new Error(\n      \"Unexpected line at position \" + lineIndex + \": \" + JSON.stringify(text) +\n      \" (expected \" + JSON.stringify(oldLines[lineIndex]) + \")\"\n    );\n  }\n\n  range.insertText(newLines[lineIndex] + trailer, \"Replace\");\n  lineIndex++;\n}\n\nawait context.sync();\n\nif (lineIndex !== oldLines.length) {\n  throw new Error(\"Only replaced \" + lineIndex + \" of \" + oldLines.length + \" lyric lines\");\n}\n", "ps1": "# Replace the lyrics of the single-paragraph song with the new\n# \"Sajeev / Nasi Goreng\" lyrics, line-by-line, while preserving the\n# original <w:br/> line-break structure inside the run.\n#\n# The whole song lives in ONE paragraph, as a single run containing many\n# text segments separated by vertical-tab line breaks (Word represents a\n# <w:br/> as Chr(11) in Range.Text). Several lines repeat verbatim\n# (e.g. \"Yaseen, Yaseen\" / \"Chorus:\"), so a plain Find/Replace would be\n# ambiguous; instead we walk the paragraph's text in document order,\n# splitting on Chr(11), and rewrite each non-blank line by its position\n# using character-offset Ranges -- applied back-to-front so earlier\n# offsets stay valid while later text lengths change.\n\n$d = $word.ActiveDocument\n\n$lines = @(\n    @{ Old = 'Verse 1:'; New = '(Verse 1)' },\n    @{ Old = 'Yaseen, my Pakistani friend'; New = 'Sajeev, my Malaysian friend,' },\n    @{ Old = 'With his warm smile that never ends'; New = 'With a love for Nasi Goreng that never ends,' },\n    @{ Old = 'From Lahore to Karachi, he''s a true gem'; New = 'He cooks it up with skill and flair,' },\n    @{ Old = 'Always there for me till the very end'; New = 'A taste of home that he loves to share.' },\n    @{ Old = 'Chorus:'; New = '(Chorus)' },\n    @{ Old = 'Yaseen, Yaseen'; New = 'Nasi Goreng, oh Nasi Goreng,' },\n    @{ Old = 'A friend so true and keen'; New = 'It''s Sajeev''s favorite thing,' },\n    @{ Old = 'In his heart, love is seen'; New = 'Spicy, savory, and oh so good,' },\n    @{ Old = 'Yaseen, Yaseen'; New = 'In every bite, you can taste the love.' },\n    @{ Old = 'Verse 2:'; New = '(Verse 2)' },\n    @{ Old = 'I remember the days we spent'; New = 'From Kuala Lumpur to Penang street,' },\n    @{ Old = 'Laughing and talking ''til the night was spent'; New = 'Sajeev knows where to find the best to eat,' },\n    @{ Old = 'His kindness and wisdom, a guiding light'; New = 'In hawker stalls and busy markets,' },\n    @{ Old = 'In his presence, everything feels right'; New = 'He''ll hunt down the dish that steals his heart.' },\n    @{ Old = 'Chorus:'; New = '(Chorus)' },\n    @{ Old = 'Yaseen, Yaseen'; New = 'Nasi Goreng, oh Nasi Goreng,' },\n    @{ Old = 'A friend so true and keen'; New = 'It''s Sajeev''s favorite thing,' },\n    @{ Old = 'In his heart, love is seen'; New = 'Spicy, savory, and oh so good,' },\n    @{ Old = 'Yaseen, Yaseen'; New = 'In every bite, you can taste the love.' },\n    @{ Old = 'Bridge:'; New = '(Bridge)' },\n    @{ Old = 'No matter the distance or time apart'; New = 'With a side of sambal and a fried egg on top,' },\n    @{ Old = 'Yaseen will always have a place in my heart'; New = 'Sajeev''s Nasi Goreng will make your taste buds pop,' },\n    @{ Old = 'With his loyalty and friendship, I am blessed'; New = 'A symphony of flavors in every spoonful,' },\n    @{ Old = 'Forever grateful for his love and zest'; New = 'It''s a dish that''s truly beautiful.' },\n    @{ Old = 'Chorus:'; New = '(Chorus)' },\n    @{ Old = 'Yaseen, Yaseen'; New = 'Nasi Goreng, oh Nasi Goreng,' },\n    @{ Old = 'A friend so true and keen'; New = 'It''s Sajeev''s favorite thing,' },\n    @{ Old = 'In his heart, love is seen'; New = 'Spicy, savory, and oh so good,' },\n    @{ Old = 'Yaseen, Yaseen'; New = 'In every bite, you can taste the love.' },\n    @{ Old = 'Outro:'; New = '(Outro)' },\n    @{ Old = 'To my dear Pakistani friend'; New = 'So here''s to Sajeev and his Nasi Goreng delight,' },\n    @{ Old = 'Yaseen, may our bond never end'; New = 'A dish that brings joy and a smile so bright,' },\n    @{ Old = 'In this world or the next, I''ll always be'; New = 'May his love for it never wane,' },\n    @{ Old = 'Grateful for your friendship eternally.'; New = 'And may he forever enjoy the taste of home again and again.' }\n)\n\n$fullText = $d.Content.Text\n$brk = [char]11\n$cr = [char]13\n\n# Content.Text ends with the story's final paragraph mark (Chr(13)); strip\n# it so the last segment's text/length lines up with the others (none of\n# which include a trailing mark).\nif ($fullText.Length -gt 0 -and $fullText[$fullText.Length - 1] -eq $cr) {\n    $fullText = $fullText.Substring(0, $fullText.Length - 1)\n}\n\n$segments = $fullText -split $brk\n\nif ($segments.Length -ne ($lines.Length + 6)) {\n    throw (\"Unexpected segment count: \" + $segments.Length)\n}\n\n# Compute the [start,end) character offset of each segment within the\n# document, then match non-blank segments (in order) against $lines.\n$offsets = New-Object System.Collections.ArrayList\n$pos = 0\nforeach ($seg in $segments) {\n    $null = $offsets.Add(@{ Start = $pos; End = ($pos + $seg.Length); Text = $seg })\n    $pos += $seg.Length + 1\n}\n\n$matches = New-Object System.Collections.ArrayList\n$li = 0\nforeach ($entry in $offsets) {\n    if ($entry.Text.Length -eq 0) {\n        continue\n    }\n    if ($li -ge $lines.Length) {\n        throw (\"More lyric lines found than expected: \" + $entry.Text)\n    }\n    if ($entry.Text -ne $lines[$li].Old) {\n        throw (\"Unexpected line at position \" + $li + \": [\" + $entry.Text + \"] (expected [\" + $lines[$li].Old + \"])\")\n    }\n    $null = $matches.Add(@{ Start = $entry.Start; End = $entry.End; New = $lines[$li].New })\n    $li++\n}\n\nif ($li -ne $lines.Length) {\n    throw (\"Only matched \" + $li + \" of \" + $lines.Length + \" lyric lines\")\n}\n\n# Apply back-to-front so that earlier offsets are not invalidated by the\n# length change of a later replacement.\nfor ($i = $matches.Count - 1; $i -ge 0; $i--) {\n    $m = $matches[$i]\n    $r = $d.Range($m.Start, $m.End)\n    $r.Text = $m.New\n}\n"}
